# Updated cryptos list on Tue Oct 31 17:43:40 UTC 2023 with GitHub Actions
#
# Refresh the crypto "Price" (column D) and "Volume(1h)" (column E) figures
# pulled from coinranking.com. The source sheet keeps these as plain text
# cells (inline strings) even when a price happens to look like a bare
# number, so for any new value that Excel would otherwise auto-convert to a
# Number on assignment, the cell's NumberFormat is switched to Text ("@")
# first to preserve the original text semantics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.399.81'
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").Value = '1.800.20'
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.19'
$ws.Range("E5").Value = '  -1.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.603'
$ws.Range("E6").Value = '  +4.37%  '
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '35.96'
$ws.Range("E8").Value = '  +3.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.291'
$ws.Range("E9").Value = '  -2.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0677'
$ws.Range("E10").Value = '  -1.99%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0966'
$ws.Range("E11").Value = '  +1.51%  '
$ws.Range("D12").Value = '2.060.77'
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.20'
$ws.Range("E13").Value = '  +0.38%  '
$ws.Range("D14").Value = '1.809.91'
$ws.Range("E14").Value = '  +0.58%  '
$ws.Range("E15").Value = '  -1.86%  '
$ws.Range("D16").Value = '34.370.13'
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.43'
$ws.Range("E17").Value = '  +2.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.50'
$ws.Range("E18").Value = '  -0.95%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.22'
$ws.Range("E19").Value = '  -0.98%  '
$ws.Range("D20").Value = '0.0₃0774'
$ws.Range("E20").Value = '  -2.81%  '
$ws.Range("E21").Value = '  -1.85%  '
$ws.Range("E22").Value = '  -0.31%  '
$ws.Range("E23").Value = '  -1.41%  '
$ws.Range("E24").Value = '  +5.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '170.17'
$ws.Range("E25").Value = '  -0.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.86'
$ws.Range("E26").Value = '  +4.41%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.35'
$ws.Range("E27").Value = '  +3.76%  '
$ws.Range("E28").Value = '  +2.48%  '
$ws.Range("E29").Value = '  -0.32%  '
$ws.Range("E30").Value = '  -1.50%  '
$ws.Range("E31").Value = '  -0.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.23'
$ws.Range("E32").Value = '  -1.51%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0516'
$ws.Range("E33").Value = '  -2.17%  '
$ws.Range("E34").Value = '  -3.13%  '
$ws.Range("D35").Value = '1.364.16'
$ws.Range("E35").Value = '  -2.49%  '
$ws.Range("E36").Value = '  -3.83%  '
$ws.Range("E37").Value = '  -0.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.37'
$ws.Range("E38").Value = '  -7.10%  '
$ws.Range("E39").Value = '  -1.46%  '
$ws.Range("E40").Value = '  +0.73%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '81.22'
$ws.Range("E41").Value = '  -1.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.78'
$ws.Range("E42").Value = '  -1.72%  '
$ws.Range("E43").Value = '  -1.25%  '
$ws.Range("E44").Value = '  +5.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.23'
$ws.Range("E45").Value = '  -3.13%  '
$ws.Range("E46").Value = '  -2.28%  '
$ws.Range("D47").Value = '1.962.53'
$ws.Range("E47").Value = '  +0.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.79'
$ws.Range("E48").Value = '  -3.41%  '
$ws.Range("E49").Value = '  -0.30%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '101.89'
$ws.Range("E50").Value = '  -2.43%  '
$ws.Range("E51").Value = '  -4.33%  '
